# Applies the "Updated cryptos list" data refresh (prices + 1h volume deltas,
# plus an InjectiveProtocol/FirstDigitalUSD row-order swap at rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value. A leading "'" forces Excel to
# keep a numeric-looking string (e.g. "0.999") as text instead of a number,
# matching the source file where every Price/Volume cell is stored as text.
$cellUpdates = [ordered]@{
    "D2" = "63.294.27"
    "E2" = "  -4.22%  "
    "D3" = "3.116.20"
    "E3" = "  -5.43%  "
    "E4" = "  -0.06%  "
    "D5" = "'559.02"
    "E5" = "  -4.79%  "
    "D6" = "'160.97"
    "E6" = "  -10.17%  "
    "E7" = "  -0.01%  "
    "D8" = "'0.580"
    "E8" = "  -9.55%  "
    "D9" = "3.108.33"
    "E9" = "  -5.67%  "
    "E10" = "  -2.40%  "
    "E11" = "  -8.39%  "
    "D12" = "'0.376"
    "E12" = "  -6.09%  "
    "D13" = "3.656.51"
    "E13" = "  -5.57%  "
    "E14" = "  -1.37%  "
    "D15" = "63.285.74"
    "E15" = "  -4.39%  "
    "D16" = "'24.59"
    "E16" = "  -7.36%  "
    "D17" = "3.104.99"
    "E17" = "  -6.32%  "
    "E18" = "  -6.82%  "
    "D19" = "'395.43"
    "E19" = "  -6.58%  "
    "D20" = "'5.18"
    "E20" = "  -5.40%  "
    "D21" = "'12.39"
    "E21" = "  -4.90%  "
    "D22" = "'6.99"
    "E22" = "  -3.96%  "
    "E23" = "  +0.25%  "
    "D24" = "'67.02"
    "E24" = "  -5.94%  "
    "D25" = "'0.199"
    "E25" = "  -4.07%  "
    "D26" = "'0.474"
    "E26" = "  -6.91%  "
    "D27" = "0.0₃0999"
    "E27" = "  -12.37%  "
    "E28" = "  -8.51%  "
    "D29" = "'0.999"
    "E29" = "  -0.19%  "
    "D30" = "'0.999"
    "E31" = "  -7.70%  "
    "D32" = "'20.78"
    "E32" = "  -6.67%  "
    "D33" = "'6.20"
    "E33" = "  -5.70%  "
    "D34" = "'4.77"
    "E34" = "  -7.51%  "
    "E35" = "  -8.18%  "
    "D36" = "'151.77"
    "E36" = "  -3.71%  "
    "E37" = "  -8.87%  "
    "D38" = "2.731.83"
    "E38" = "  -4.69%  "
    "E39" = "  -8.93%  "
    "D40" = "'4.04"
    "E40" = "  -6.77%  "
    "D41" = "'23.18"
    "E41" = "  -11.94%  "
    "D42" = "'38.19"
    "E42" = "  -3.67%  "
    "D43" = "'0.691"
    "E43" = "  -7.75%  "
    "D44" = "'0.0608"
    "E44" = "  -4.67%  "
    "E45" = "  -7.87%  "
    "E46" = "  -5.06%  "
    "B47" = "InjectiveProtocol"
    "C47" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D47" = "'20.66"
    "E47" = "  -9.73%  "
    "B48" = "FirstDigitalUSD"
    "C48" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D48" = "'0.999"
    "E48" = "  -0.06%  "
    "D49" = "'278.14"
    "E49" = "  -11.11%  "
    "D50" = "'0.0971"
    "E50" = "  -4.74%  "
    "E51" = "  +0.48%  "
}

foreach ($ref in $cellUpdates.Keys) {
    $ws.Range($ref).Value = $cellUpdates[$ref]
}
